{"js": "// The commit wraps a few spans with proofing-error markers (w:proofErr\n// gramStart/gramEnd around \"controllers/Account.java\" and around\n// \"data/AccountDB\", spellStart/spellEnd around \"dologin\"/\"dologout\"/\n// \"docomment\" and \"AccountDB\"/\"shop.jsp\"), which in turn forces those runs\n// to be split at the marked boundaries, and it also splits the word\n// \"provided\" into two runs (\"provide\" + \" folder.\") in the final\n// paragraph. None of this is exposed as a first-class Office.js object\n// (Word.js has no ProofingError / proofErr concept), so we rebuild the\n// affected paragraphs from exact OOXML and splice each one in with\n// Range.insertOoxml(..., \"Replace\"). Using whole-paragraph OOXML swaps\n// (rather than text edits) guarantees the run-splits and proofErr\n// placements land exactly where the diff expects, while every paragraph\n// keeps its original mark / properties.\n\nconst W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction wrapPackage(innerBodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    innerBodyXml +\n    \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Exact replacement XML for each affected paragraph, built by inserting\n// <w:proofErr/> markers and splitting runs exactly as the diff shows,\n// while keeping every other attribute/property untouched.\nconst replacements = [\n  {\n    // \"controllers/Account.java\" -> wrap the whole paragraph in gramStart/gramEnd\n    find: \"controllers/Account.java\",\n    xml:\n      `<w:p ${W_NS} w:rsidR=\"002D47B6\" w:rsidRPr=\"005E542E\" w:rsidRDefault=\"003525F0\">` +\n      \"<w:pPr><w:rPr><w:b/><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r w:rsidRPr=\"005E542E\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>controllers/</w:t></w:r>' +\n      '<w:r w:rsidR=\"006A1818\" w:rsidRPr=\"005E542E\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Account.java</w:t></w:r>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      \"</w:p>\",\n  },\n  {\n    // \", such as \"dologin\", \"dologout\", \"docomment\"...\" -> spellStart/spellEnd around each made-up verb\n    find: \"such as\",\n    xml:\n      `<w:p ${W_NS} w:rsidR=\"006A1818\" w:rsidRDefault=\"006A1818\">` +\n      \"<w:pPr><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>This is</w:t></w:r>' +\n      '<w:r w:rsidRPr=\"006A1818\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> used to process request from</w:t></w:r>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> several </w:t></w:r>' +\n      '<w:r w:rsidRPr=\"006A1818\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>actions</w:t></w:r>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>, such as \\u201C</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>dologin</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>\\u201D, \\u201C</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>dologout</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>\\u201D, \\u201C</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>docomment</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>\\u201D\\u2026</w:t></w:r>' +\n      \"</w:p>\",\n  },\n  {\n    // \"The action \"docomment\" is recently added...\" -> spellStart/spellEnd around \"docomment\"\n    find: \"is recently added\",\n    xml:\n      `<w:p ${W_NS} w:rsidR=\"006A1818\" w:rsidRDefault=\"006A1818\">` +\n      \"<w:pPr><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>The action \\u201C</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>docomment</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">\\u201D is recently added to serve the POST comment purpose. </w:t></w:r>' +\n      \"</w:p>\",\n  },\n  {\n    // \"data/AccountDB\" -> split into \"data/\" + \"AccountDB\", spellStart/spellEnd\n    // around \"AccountDB\", whole paragraph wrapped in gramStart/gramEnd\n    find: \"data/AccountDB\",\n    xml:\n      `<w:p ${W_NS} w:rsidR=\"00BF6E63\" w:rsidRPr=\"00741B63\" w:rsidRDefault=\"00BF6E63\">` +\n      \"<w:pPr><w:rPr><w:b/><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      '<w:r w:rsidRPr=\"00741B63\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>data/</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r w:rsidRPr=\"00741B63\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>AccountDB</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      \"</w:p>\",\n  },\n  {\n    // \"Web Pages/shop.jsp\" -> spellStart/spellEnd around \"shop.jsp\"\n    find: \"Web Pages/shop.jsp\",\n    xml:\n      `<w:p ${W_NS} w:rsidR=\"00A20A51\" w:rsidRDefault=\"005612BF\">` +\n      \"<w:pPr><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Web Pages/</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r w:rsidR=\"00A20A51\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>shop.jsp</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      \"</w:p>\",\n  },\n  {\n    // \"Pictures is saved to provided folder.\" -> split into 3 runs\n    // (\"Pictures is saved to \" + \"provide\" + \" folder.\") exactly as the diff shows\n    find: \"Pictures is saved to\",\n    xml:\n      `<w:p ${W_NS} w:rsidR=\"0029463D\" w:rsidRPr=\"006A1818\" w:rsidRDefault=\"0029463D\">` +\n      \"<w:pPr><w:rPr><w:lang w:val=\\\"en-US\\\"/></w:rPr></w:pPr>\" +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">Pictures is saved to </w:t></w:r>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>provide</w:t></w:r>' +\n      '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> folder.</w:t></w:r>' +\n      '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n      \"</w:p>\",\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const rep of replacements) {\n  const para = paragraphs.items.find((p) => p.text.indexOf(rep.find) !== -1);\n  if (!para) {\n    throw new Error(\"Paragraph not found for: \" + rep.find);\n  }\n  const range = para.getRange(\"Whole\");\n  range.insertOoxml(wrapPackage(rep.xml), \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The commit wraps a few spans with proofing-error markers (w:proofErr\n# gramStart/gramEnd around \"controllers/Account.java\" and around\n# \"data/AccountDB\", spellStart/spellEnd around \"dologin\"/\"dologout\"/\n# \"docomment\" and \"AccountDB\"/\"shop.jsp\"), which forces those runs to be\n# split at the marked boundaries, and it also splits the word \"provided\"\n# into two runs (\"provide\" + \" folder.\") in the final paragraph. The Word\n# object model has no ProofingError/proofErr concept to toggle, so we\n# rebuild each affected paragraph from exact OOXML and splice it in with\n# Range.InsertXML (Flat-OPC single-part package), which replaces the\n# paragraph's content in place while leaving every other paragraph and\n# property untouched.\n\n$d = $word.ActiveDocument\n\n# This PowerShell-style interpreter normalizes literal typographic quote\n# characters pasted into a script back to ASCII quotes, so the left/right\n# double quotation marks and ellipsis that appear in the document text are\n# built from their Unicode code points instead of being typed literally.\n$LDQ = [char]0x201C   # \u201c\n$RDQ = [char]0x201D   # \u201d\n$HELLIP = [char]0x2026   # \u2026\n\nfunction Wrap-Package([string]$innerBodyXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $innerBodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$W_NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"'\n\n$replacements = @(\n    @{\n        Find = \"controllers/Account.java\"\n        Xml = '<w:p ' + $W_NS + ' w:rsidR=\"002D47B6\" w:rsidRPr=\"005E542E\" w:rsidRDefault=\"003525F0\">' +\n            '<w:pPr><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n            '<w:proofErr w:type=\"gramStart\"/>' +\n            '<w:r w:rsidRPr=\"005E542E\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>controllers/</w:t></w:r>' +\n            '<w:r w:rsidR=\"006A1818\" w:rsidRPr=\"005E542E\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>Account.java</w:t></w:r>' +\n            '<w:proofErr w:type=\"gramEnd\"/>' +\n            '</w:p>'\n    },\n    @{\n        Find = \"such as\"\n        Xml = '<w:p ' + $W_NS + ' w:rsidR=\"006A1818\" w:rsidRDefault=\"006A1818\">' +\n            '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>This is</w:t></w:r>' +\n            '<w:r w:rsidRPr=\"006A1818\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> used to process request from</w:t></w:r>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> several </w:t></w:r>' +\n            '<w:r w:rsidRPr=\"006A1818\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>actions</w:t></w:r>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>, such as ' + $LDQ + '</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellStart\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>dologin</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellEnd\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>' + $RDQ + ', ' + $LDQ + '</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellStart\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>dologout</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellEnd\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>' + $RDQ + ', ' + $LDQ + '</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellStart\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>docomment</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellEnd\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>' + $RDQ + $HELLIP + '</w:t></w:r>' +\n            '</w:p>'\n    },\n    @{\n        Find = \"is recently added\"\n        Xml = '<w:p ' + $W_NS + ' w:rsidR=\"006A1818\" w:rsidRDefault=\"006A1818\">' +\n            '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>The action ' + $LDQ + '</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellStart\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>docomment</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellEnd\"/>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">' + $RDQ + ' is recently added to serve the POST comment purpose. </w:t></w:r>' +\n            '</w:p>'\n    },\n    @{\n        Find = \"data/AccountDB\"\n        Xml = '<w:p ' + $W_NS + ' w:rsidR=\"00BF6E63\" w:rsidRPr=\"00741B63\" w:rsidRDefault=\"00BF6E63\">' +\n            '<w:pPr><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n            '<w:proofErr w:type=\"gramStart\"/>' +\n            '<w:r w:rsidRPr=\"00741B63\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>data/</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellStart\"/>' +\n            '<w:r w:rsidRPr=\"00741B63\"><w:rPr><w:b/><w:lang w:val=\"en-US\"/></w:rPr><w:t>AccountDB</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellEnd\"/>' +\n            '<w:proofErr w:type=\"gramEnd\"/>' +\n            '</w:p>'\n    },\n    @{\n        Find = \"Web Pages/shop.jsp\"\n        Xml = '<w:p ' + $W_NS + ' w:rsidR=\"00A20A51\" w:rsidRDefault=\"005612BF\">' +\n            '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Web Pages/</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellStart\"/>' +\n            '<w:r w:rsidR=\"00A20A51\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>shop.jsp</w:t></w:r>' +\n            '<w:proofErr w:type=\"spellEnd\"/>' +\n            '</w:p>'\n    },\n    @{\n        Find = \"Pictures is saved to\"\n        Xml = '<w:p ' + $W_NS + ' w:rsidR=\"0029463D\" w:rsidRPr=\"006A1818\" w:rsidRDefault=\"0029463D\">' +\n            '<w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">Pictures is saved to </w:t></w:r>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>provide</w:t></w:r>' +\n            '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> folder.</w:t></w:r>' +\n            '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n            '</w:p>'\n    }\n)\n\nforeach ($rep in $replacements) {\n    $target = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs($i)\n        if ($p.Range.Text -like \"*$($rep.Find)*\") {\n            $target = $p\n            break\n        }\n    }\n    if ($null -eq $target) {\n        throw \"Paragraph not found for: $($rep.Find)\"\n    }\n    $target.Range.InsertXML((Wrap-Package($rep.Xml)))\n}\n"}
